$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Saturday's time-in / time-out entries (row 18): 1:00 PM in, 1:30 PM out.
$ws.Range("C18").Value = 0.541666666666667
$ws.Range("D18").Value = 0.5625

# Move the active selection to E18, matching the recorded cursor position
# after the entries were made.
$ws.Range("E18").Select()
